$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.232.24"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "1.674.60"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("E4").Value = "  -0.63%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5266"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.88%  "
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("E8").Value = "  -3.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06279"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.69%  "
$ws.Range("E10").Value = "  -3.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07568"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("D12").Value = "1.677.81"
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.464"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5611"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008005"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.83%  "
$ws.Range("D17").Value = "26.038.66"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.816"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "187.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.94%  "
$ws.Range("E21").Value = "  -5.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.212"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("E23").Value = "  -0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1253"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06223"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.360"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.507"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.01%  "
$ws.Range("E32").Value = "  -4.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.634"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("E34").Value = "  -3.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6059"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.755"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.115"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01619"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("D40").Value = "1.101.91"
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8712"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.006"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("D44").Value = "1.826.93"
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000111"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.010"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05234"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.981"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.29%  "
